# Update capital structure database values for Mauritius / Precious Metals rows (2 and 3).
# Both data rows receive identical updates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2, 3) {
    # trailing_net_income
    $ws.Cells.Item($row, 11).Value = -1.42          # K

    # cash
    $ws.Cells.Item($row, 21).Value = 0.348           # U
    # cash_market_cap
    $ws.Cells.Item($row, 22).Value = 0.04438775510204081   # V

    # roe - removed in new data
    $ws.Cells.Item($row, 23).ClearContents()         # W

    # cost_equity
    $ws.Cells.Item($row, 24).Value = 0.06361838375065348   # X

    # roe_cost_equity - removed in new data
    $ws.Cells.Item($row, 25).ClearContents()         # Y
    # sales_invested_capital - removed in new data
    $ws.Cells.Item($row, 26).ClearContents()         # Z
    # roic - removed in new data
    $ws.Cells.Item($row, 27).ClearContents()         # AA

    # cost_capital
    $ws.Cells.Item($row, 28).Value = 0.06361838375065348   # AB

    # roic_cost_capital - removed in new data
    $ws.Cells.Item($row, 29).ClearContents()         # AC

    # net_debt
    $ws.Cells.Item($row, 33).Value = -0.348          # AG

    # debt_book_capital
    $ws.Cells.Item($row, 35).Value = -0              # AI

    # net_debt_market_capital
    $ws.Cells.Item($row, 36).Value = -0.04644954618259477  # AJ

    # net_debt_book_capital
    $ws.Cells.Item($row, 37).Value = 0.5631067961165048    # AK

    # interest_expenses
    $ws.Cells.Item($row, 38).Value = 0.056           # AL

    # net_interest_expenses
    $ws.Cells.Item($row, 39).Value = 0.025           # AM

    # ebit_interest_expenses - newly added value
    $ws.Cells.Item($row, 41).Value = -10.66071428571428    # AO

    # ebit_net_interest_expenses
    $ws.Cells.Item($row, 43).Value = -23.88          # AQ
}
